# Updates the cryptos price/volume table with refreshed values from the
# upstream scrape. Price/volume cells hold number-like text (e.g. "67.019.59",
# "  -1.15%  ") that must stay text, so each is written with a leading
# apostrophe (forces text entry, same as typing '1.00 into Excel) and then
# the cell style is reset to "Normal" so no stray quote-prefix/number-format
# style sticks around on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''67.019.59'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -1.15%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''2.471.73'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -2.28%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  -0.05%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''582.60'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -1.31%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''169.00'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -2.45%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  +0.02%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''0.513'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  -1.96%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''2.472.38'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -2.20%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = '''  -2.37%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = '''  -0.77%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''4.92'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -2.11%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = '''  -3.91%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '''2.924.87'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -3.09%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").Value = '''25.47'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -3.35%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''66.889.09'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -1.32%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = '''  -4.09%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''2.465.40'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -2.38%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''11.04'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -6.14%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''7.51'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -6.14%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''350.61'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -5.16%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''4.04'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -2.33%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  -0.01%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''68.83'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -4.17%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''4.24'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  -6.94%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''1.82'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -4.93%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''9.17'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -7.83%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = '''  -58.90%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D30").Value = '''0.0₃0901'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -6.49%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''508.15'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -5.52%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''7.69'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  -7.28%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").Value = '''1.76'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  -5.18%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").Value = '''1.23'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  -6.72%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  +0.03%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''159.40'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  -0.43%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''0.115'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  -10.39%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = '''  +0.27%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''18.29'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  -4.64%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = '''  -7.61%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = '''  -4.72%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = '''  +0.09%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''4.80'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  -6.34%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''0.327'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  -6.52%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''2.36'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -6.22%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = '''  -0.97%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''140.76'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -4.71%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''3.44'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -7.16%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''0.512'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  -6.69%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = '''  -10.44%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''1.58'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  -7.23%  '
$ws.Range("E51").Style = "Normal"
